$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, pushing existing rows 47..73 down to 48..74
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new record's data.
# (Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across this sheet's dataset.)
$ws.Cells.Item(47, 1).Value = 10
$ws.Cells.Item(47, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(47, 3).Value = "La Araucanía"
$ws.Cells.Item(47, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(47, 5).Value = 9
$ws.Cells.Item(47, 6).Value = 100112012
$ws.Cells.Item(47, 7).Value = "Espinaca"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 30
$ws.Cells.Item(47, 11).Value = 9000
$ws.Cells.Item(47, 12).Value = 9000
$ws.Cells.Item(47, 13).Value = 9000
$ws.Cells.Item(47, 14).Value = "$/docena de atados"
$ws.Cells.Item(47, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(47, 16).Value = 3000
$ws.Cells.Item(47, 17).Value = 3
$ws.Cells.Item(47, 18).Value = "Hortaliza"
